$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("прямая")

# Fixed error in F: correct the diagonal error-norm values in the first block.
$ws.Range("C4").Value = [double]"2.9743052456139E-9"
$ws.Range("D5").Value = [double]"7.9886832755140596E-10"
$ws.Range("E6").Value = [double]"2.80243981817775E-10"
$ws.Range("F7").Value = [double]"1.34755068097749E-10"
$ws.Range("G8").Value = [double]"6.74674824515285E-11"

# Update the active cell selection on the sheet to C4
[void]$ws.Range("C4").Select()
